# setup gui with tkinter
# Update the VOLVO compatibility table: drop the C40 and XC40 rows,
# refresh year ranges / positions / engine notes for the remaining models.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table shrinks from 9 data rows (C40, S60, S90, V60, V90,
# V90 CROSS COUNTRY, XC40, XC60, XC90) down to 7 data rows (S60, S90, V60,
# V90, V90 CROSS COUNTRY, XC60, XC90), so remove the two trailing rows first.
$ws.Rows("9:10").Delete()

# Row 2: was C40 / 2022-2023 / Rear or Front -> now S60 / 2019-2020 / Front
$ws.Range("A2").Value = "VOLVO"
$ws.Range("B2").Value = "S60"
$ws.Range("C2").Value = "2019-2020"
$ws.Range("D2").Value = "Front"
$ws.Range("E2").ClearContents()

# Row 3: was S60 / 2019-2023 / Rear or Front / AWD -> now S90 / 2017-2019 / Front
$ws.Range("A3").Value = "VOLVO"
$ws.Range("B3").Value = "S90"
$ws.Range("C3").Value = "2017-2019"
$ws.Range("D3").Value = "Front"
$ws.Range("E3").ClearContents()

# Row 4: was S90 / 2017-2023 / Rear or Front -> now V60 / 2019-2020, no position, with engine note
$ws.Range("A4").Value = "VOLVO"
$ws.Range("B4").Value = "V60"
$ws.Range("C4").Value = "2019-2020"
$ws.Range("D4").ClearContents()
$ws.Range("E4").Value = "No 2.0L L4 ELECTRIC/GAS Turbocharged"

# Row 5: was V60 / 2019-2023 / Rear or Front -> now V90 / 2018-2019, no position/engine
$ws.Range("A5").Value = "VOLVO"
$ws.Range("B5").Value = "V90"
$ws.Range("C5").Value = "2018-2019"
$ws.Range("D5").ClearContents()
$ws.Range("E5").ClearContents()

# Row 6: was V90 / 2018-2023 / Rear or Front -> now V90 CROSS COUNTRY / 2017 / Front
$ws.Range("A6").Value = "VOLVO"
$ws.Range("B6").Value = "V90 CROSS COUNTRY"
$ws.Range("C6").Value = "2017"
$ws.Range("D6").Value = "Front"
$ws.Range("E6").ClearContents()

# Row 7: was V90 CROSS COUNTRY / 2017-2021 / Rear or Front -> now XC60 / 2018-2020 / Front
$ws.Range("A7").Value = "VOLVO"
$ws.Range("B7").Value = "XC60"
$ws.Range("C7").Value = "2018-2020"
$ws.Range("D7").Value = "Front"
$ws.Range("E7").ClearContents()

# Row 8: was XC40 / 2021-2023 / Rear or Front / No 2.0L L4 Turbocharged -> now XC90 / 2016-2019 / Front
$ws.Range("A8").Value = "VOLVO"
$ws.Range("B8").Value = "XC90"
$ws.Range("C8").Value = "2016-2019"
$ws.Range("D8").Value = "Front"
$ws.Range("E8").ClearContents()
